# (JMT) Added coverage for bl_1s11 and bl_1s13
# Update the user-specific workspace/config paths used by this block
# (from user4's paths to user2's paths), and move the active selection
# from F2 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update proj_path / config_path values in row 3 (order matters for
# shared-string table placement: these two first, then the template path).
$ws.Range("J3").Value = "/pub/home/user2/jmt_workspace"
$ws.Range("K3").Value = "/pub/home/user2/jmt_workspace/workshop_config.sdl"

# Update the template path shown in B1.
$ws.Range("B1").Value = "template /pub/home/user2/jmt_workspace/blocks/bl_1s11/bl_1s11.tsdl"

# Move the selected/active cell to B2.
[void]$ws.Range("B2").Select()
